$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "35.157.29"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.31%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.854.07"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  +0.56%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "237.91"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.13%  "
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  +0.52%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "41.78"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +4.94%  "
$ws.Range("E9").Value = "  +1.10%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0691"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.02%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0989"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.19%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "2.122.15"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.83%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.859.17"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("E16").Value = "  +1.46%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "35.140.85"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.27%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "69.94"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.34%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.0₃0791"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.43%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "240.53"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.23%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.15"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("E24").Value = "  +0.63%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "169.06"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.61%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.96"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.94%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.83"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +20.63%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "17.55"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.94%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.123"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("E30").Value = "  +0.57%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0554"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.57%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.97"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.65%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.00"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.99%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.75"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +28.50%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.99"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +8.74%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.802"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +15.32%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.30"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +3.46%  "
$ws.Range("E38").Value = "  +8.16%  "
$ws.Range("E39").Value = "  +3.21%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "89.69"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -3.29%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.341.56"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.01%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "14.82"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.45%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "12.87"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +47.32%  "
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("E45").Value = "  +0.68%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0553"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +6.50%  "
$ws.Range("E47").Value = "  -0.73%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "6.44"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.20%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.038.81"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("E51").Value = "  +1.56%  "
